# Add a switching capacitor (Murata GRM31C5C1E104JA01K, 1206 package) to the
# part list, replacing the previous 0603-package capacitor entries for C5/C6,
# and add the Mouser order number + product link for the new part.
#
# Shared-string cells are written in the same order the new unique strings
# appear in the target workbook (link, order-number, package-name,
# device-name) so the newly minted shared-string table entries line up with
# the diff: 149=link, 150=order number, 151=C1206, 152=C-EUC1206 Zwitschi-Cap.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$link = "https://www.mouser.de/ProductDetail/Murata-Electronics/GRM31C5C1E104JA01K?qs=%2Fha2pyFaduhBfEIR6jcOaJ56CU6eddb5oLEFH%252BIY%2FeBTD895l24yhA%3D%3D"
$orderNo = "81-GRM31C5C1E104JA1K "
$package = "C1206"
$device = "C-EUC1206 Zwitschi-Cap"

# Row 15 = part C5
$ws.Range("I15").Value = $link
$ws.Range("H15").Value = $orderNo
$ws.Range("E15").Value = $package
$ws.Range("D15").Value = $device

# Row 16 = part C6
$ws.Range("I16").Value = $link
$ws.Range("H16").Value = $orderNo
$ws.Range("E16").Value = $package
$ws.Range("D16").Value = $device

# Best-effort view-state update to mirror the author's selection/scroll.
$ws.Activate() | Out-Null
$ws.Range("F17").Select() | Out-Null
try {
    $excel.ActiveWindow.ScrollRow = 4
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
    # scroll position isn't modelled in all hosts - ignore if unsupported
}
